$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 393, shifting existing rows 393-398 down to 395-400
$ws.Rows("393:394").Insert()

# New row 393 data
$ws.Range("A393").Value = 3
$ws.Range("B393").Value = "Femacal de La Calera"
$ws.Range("C393").Value = "Coquimbo"
$ws.Range("D393").Value = 44628
$ws.Range("E393").Value = 5
$ws.Range("F393").Value = 100112032
$ws.Range("G393").Value = "Zapallo italiano"
$ws.Range("H393").Value = "Sin especificar"
$ws.Range("I393").Value = "Primera"
$ws.Range("J393").Value = 125
$ws.Range("K393").Value = 5000
$ws.Range("L393").Value = 5500
$ws.Range("M393").Value = 5260
$ws.Range("N393").Value = "$/caja 36 unidades"
$ws.Range("O393").Value = "Provincia de Quillota"
$ws.Range("P393").Value = 146
$ws.Range("Q393").Value = 36
$ws.Range("R393").Value = "Hortaliza"

# New row 394 data
$ws.Range("A394").Value = 3
$ws.Range("B394").Value = "Femacal de La Calera"
$ws.Range("C394").Value = "Coquimbo"
$ws.Range("D394").Value = 44628
$ws.Range("E394").Value = 5
$ws.Range("F394").Value = 100112032
$ws.Range("G394").Value = "Zapallo italiano"
$ws.Range("H394").Value = "Sin especificar"
$ws.Range("I394").Value = "Primera"
$ws.Range("J394").Value = 165
$ws.Range("K394").Value = 9000
$ws.Range("L394").Value = 10000
$ws.Range("M394").Value = 9515
$ws.Range("N394").Value = "$/caja 70 unidades"
$ws.Range("O394").Value = "Provincia de Quillota"
$ws.Range("P394").Value = 136
$ws.Range("Q394").Value = 70
$ws.Range("R394").Value = "Hortaliza"
